$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.462.33"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "1.839.34"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'260.33"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "'0.5254"
$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("D8").Value = "'0.3191"
$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("D9").Value = "'0.06788"
$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("D10").Value = "'18.75"
$ws.Range("E10").Value = "  +1.12%  "

$ws.Range("D11").Value = "'0.7854"
$ws.Range("E11").Value = "  +2.94%  "

$ws.Range("D12").Value = "'0.07744"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").Value = "1.831.82"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").Value = "'87.71"
$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("D15").Value = "'5.013"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Value = "'13.83"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D19").Value = "'0.000007949"
$ws.Range("E19").Value = "  +0.44%  "

$ws.Range("D20").Value = "26.490.93"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "2.070.93"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'4.624"
$ws.Range("E22").Value = "  +1.55%  "

$ws.Range("D23").Value = "'5.975"
$ws.Range("E23").Value = "  +0.75%  "

$ws.Range("D24").Value = "'9.369"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("E25").Value = "  -2.53%  "

$ws.Range("D26").Value = "'2.177"
$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").Value = "'16.92"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "'111.46"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").Value = "'4.155"
$ws.Range("E30").Value = "  -0.35%  "

$ws.Range("D31").Value = "'0.08687"
$ws.Range("E31").Value = "  -0.52%  "

$ws.Range("E32").Value = "  -1.46%  "

$ws.Range("D33").Value = "'0.04874"
$ws.Range("E33").Value = "  +1.26%  "

$ws.Range("D34").Value = "'0.7280"
$ws.Range("E34").Value = "  +4.20%  "

$ws.Range("D35").Value = "'1.134"
$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("D36").Value = "'2.861"
$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("D38").Value = "'2.246"
$ws.Range("E38").Value = "  +2.39%  "

$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("D40").Value = "'0.4776"
$ws.Range("E40").Value = "  -1.15%  "

$ws.Range("D41").Value = "'0.8927"
$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("D42").Value = "'109.52"
$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("E43").Value = "  -2.59%  "

$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").Value = "'7.669"
$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("D46").Value = "'0.4169"
$ws.Range("E46").Value = "  +1.55%  "

$ws.Range("D47").Value = "'8.986"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").Value = "'0.05850"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "'0.1231"
$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("D50").Value = "'34.84"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").Value = "'0.8907"
$ws.Range("E51").Value = "  +1.33%  "
